# Rerun all TODE norms, found error in grade code that was creating the wierd column
# Apply corrected scaled-score (ss) values to the raw-score lookup tables on each age tab.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { cell -> new value }
$changes = @{
    "5.0-5.3" = @{
        "B5"  = 91
        "B9"  = 98
        "B10" = 100
        "B11" = 102
        "B12" = 104
        "B13" = 106
        "B14" = 108
        "B15" = 110
        "B22" = 129
    }
    "5.4-5.7" = @{
        "B8"  = 91
        "B12" = 98
        "B13" = 100
        "B17" = 109
        "B21" = 121
        "B22" = 124
        "B23" = 127
    }
    "5.8-5.11" = @{
        "B8"  = 86
        "B12" = 93
        "B13" = 95
        "B14" = 97
        "B22" = 119
        "B24" = 126
    }
    "6.0-6.5" = @{
        "B2"  = 73
        "B7"  = 79
        "B10" = 84
        "B21" = 110
        "B23" = 117
    }
    "6.6-6.11" = @{
        "B3"  = 69
        "B10" = 78
        "B13" = 83
        "B14" = 85
        "B15" = 87
        "B22" = 107
        "B26" = 129
    }
    "7.0-7.5" = @{
        "B2"  = 63
        "B12" = 76
        "B18" = 89
        "B24" = 111
        "B25" = 118
        "B26" = 127
    }
    "7.6-7.11" = @{
        "B2"  = 59
        "B20" = 89
        "B24" = 108
        "B25" = 116
        "B26" = 127
    }
    "8.0-8.5" = @{
        "B5"  = 59
        "B15" = 73
        "B16" = 75
        "B17" = 77
        "B18" = 79
        "B21" = 88
        "B22" = 92
        "B24" = 104
        "B25" = 115
        "B26" = 126
        "B27" = 130
    }
    "8.6-9.3" = @{
        "B8"  = 58
        "B19" = 76
        "B20" = 79
        "B21" = 82
        "B23" = 92
        "B24" = 102
        "B25" = 115
        "B26" = 123
        "B27" = 128
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $changes[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
